$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.997.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.540.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.52%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.13%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.215"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.664"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000308"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.110.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "615.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.20%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.061.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.558.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.49%  "

$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("E21").Value = "  -1.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.65%  "

$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "538.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.400"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.60%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.33%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0781"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.550.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("E43").Value = "  +3.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0464"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.86%  "

$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("E46").Value = "  +4.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.55%  "

$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.58%  "
